# Apply the edit described by the diff:
#  - Insert a new data row for "FEMAPENT 2.5MG 14 F.C. TABLETS" right before
#    the existing "FLACORT 30MG 20 TAB" row (which is row 11 in the sheet).
#  - All rows below shift down by one (Excel handles this automatically when
#    a whole row is inserted).
#  - Renumber the running sequence number (column A) for the new row and
#    refresh the grand-total cell that sits a couple of rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank row at row 11 (this shifts rows 11..34 down to 12..35,
#    including merged cells, the totals row and the footer row).
$ws.Rows("11:11").Insert()

# 2. Give the freshly inserted row the same merged layout as every other
#    data row (B:G, H:K and L:M are merged on every product row).
$ws.Range("B11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()

# 3. Copy the cell formatting (styles/borders/fonts) from the row that is
#    now directly below (the old row 11, now row 12) so the new row looks
#    identical to the rest of the table instead of using blank default
#    formatting.
$ws.Range("A12:N12").Copy()
$ws.Range("A11:N11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Match the row height used by the surrounding rows.
$ws.Rows("11:11").RowHeight = 24.75

# 5. Fill in the values for the new product row.
$ws.Cells.Item(11, 1).Value2 = 8
$ws.Cells.Item(11, 2).Value2 = "FEMAPENT 2.5MG 14 F.C. TABLETS"
$ws.Cells.Item(11, 8).Value2 = "1:1"
$ws.Cells.Item(11, 12).Value2 = 197
$ws.Cells.Item(11, 14).Value2 = "0:0"

# 6. Renumber column A (the "م" sequence column) for every row that was
#    pushed down, so it keeps counting 1,2,3... without a gap/duplicate.
for ($r = 12; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 3
}

# 7. Update the grand-total cell (K column on the totals row, now row 33)
#    to include the new row's amount (was 1219.87, now +197 = 1416.87).
$ws.Cells.Item(33, 11).Value2 = 1416.8699999999999
